# Updates a handful of numeric cells in the "Resumo Inscricoes Integrado"
# sheet (Inscritos / Pagos / Inscricoes homologadas counts) to reflect the
# latest registration tallies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 95
$ws.Range("F2").Value = 65
$ws.Range("H2").Value = 71

$ws.Range("E10").Value = 522
$ws.Range("E11").Value = 339
$ws.Range("E12").Value = 513

$ws.Range("F14").Value = 65
$ws.Range("H14").Value = 100

$ws.Range("E15").Value = 161
$ws.Range("F15").Value = 71
$ws.Range("H15").Value = 120

$ws.Range("E31").Value = 72

$ws.Range("E33").Value = 288
$ws.Range("E34").Value = 212

$ws.Range("E37").Value = 155

$ws.Range("E40").Value = 256
$ws.Range("E41").Value = 381
$ws.Range("E42").Value = 374
$ws.Range("E43").Value = 114

$ws.Range("E45").Value = 141

$ws.Range("E47").Value = 445
$ws.Range("F47").Value = 224
$ws.Range("H47").Value = 316

$ws.Range("E48").Value = 203
$ws.Range("E49").Value = 283

$ws.Range("E50").Value = 239
$ws.Range("F50").Value = 114
$ws.Range("H50").Value = 187

$ws.Range("F51").Value = 100
$ws.Range("H51").Value = 172
